$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69
$ws.Range("B69").Value = 7423700
$ws.Range("E69").Value = "TuS Hornau"
$ws.Range("F69").Value = "FC Burgsolms"
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 1
$ws.Range("K69").Value = "H"
$ws.Range("L69").Value = 1.727
$ws.Range("M69").Value = 4.5
$ws.Range("N69").Value = 3.2
$ws.Range("O69").Value = 1.727
$ws.Range("P69").Value = 4.5
$ws.Range("Q69").Value = 3.2
$ws.Range("R69").Value = -0.5
$ws.Range("S69").Value = 1.775
$ws.Range("T69").Value = 2.025
$ws.Range("U69").Value = 3.5
$ws.Range("V69").Value = 1.85
$ws.Range("W69").Value = 1.95
$ws.Range("X69").Value = 0.7270000000000001
$ws.Range("Y69").Value = -1
$ws.Range("AA69").Value = 0.7749999999999999
$ws.Range("AB69").Value = -1
$ws.Range("AC69").Value = -1
$ws.Range("AD69").Value = 0.95

# Row 70
$ws.Range("B70").Value = 7423699
$ws.Range("E70").Value = "SG 2000 MulheimKarlich"
$ws.Range("F70").Value = "Ahrweiler BC"
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 2
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = "D"
$ws.Range("L70").Value = 2.2
$ws.Range("M70").Value = 5
$ws.Range("N70").Value = 2.2
$ws.Range("O70").Value = 2.2
$ws.Range("P70").Value = 4.75
$ws.Range("Q70").Value = 2.2
$ws.Range("R70").Value = 0
$ws.Range("S70").Value = 1.9
$ws.Range("T70").Value = 1.9
$ws.Range("U70").Value = 4.25
$ws.Range("V70").Value = 1.775
$ws.Range("W70").Value = 2.025
$ws.Range("X70").Value = -1
$ws.Range("Y70").Value = 3.75
$ws.Range("AA70").Value = 0
$ws.Range("AB70").Value = 0
$ws.Range("AC70").Value = -0.5
$ws.Range("AD70").Value = 0.5125

# Row 125
$ws.Range("B125").Value = 8039381
$ws.Range("F125").Value = "SV 1908 GW Ahrensfelde"
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 3
$ws.Range("I125").Value = 1
$ws.Range("J125").Value = 3
$ws.Range("L125").Value = 3.25
$ws.Range("M125").Value = 3.8
$ws.Range("N125").Value = 1.833
$ws.Range("O125").Value = 3.25
$ws.Range("P125").Value = 3.8
$ws.Range("Q125").Value = 1.833
$ws.Range("R125").Value = 0.5
$ws.Range("S125").Value = 1.925
$ws.Range("T125").Value = 1.875
$ws.Range("U125").Value = 3
$ws.Range("V125").Value = 1.825
$ws.Range("W125").Value = 1.975
$ws.Range("Z125").Value = 0.833
$ws.Range("AB125").Value = 0.875
$ws.Range("AC125").Value = 0.825
$ws.Range("AD125").Value = -1

# Row 126
$ws.Range("B126").Value = 8039382
$ws.Range("F126").Value = "TSV Steinbach II"
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 4
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2
$ws.Range("L126").Value = 15
$ws.Range("M126").Value = 9
$ws.Range("N126").Value = 1.111
$ws.Range("O126").Value = 15
$ws.Range("P126").Value = 9
$ws.Range("Q126").Value = 1.111
$ws.Range("R126").Value = 2.75
$ws.Range("S126").Value = 1.9
$ws.Range("T126").Value = 1.9
$ws.Range("U126").Value = 4
$ws.Range("V126").Value = 1.9
$ws.Range("W126").Value = 1.9
$ws.Range("Z126").Value = 0.111
$ws.Range("AB126").Value = 0.8999999999999999
$ws.Range("AC126").Value = 0
$ws.Range("AD126").Value = 0

# Row 180
$ws.Range("B180").Value = 8271388
$ws.Range("F180").Value = "SSC Weissenfels"
$ws.Range("G180").Value = 4
$ws.Range("L180").Value = 3.4
$ws.Range("M180").Value = 4.2
$ws.Range("N180").Value = 1.727
$ws.Range("O180").Value = 3.4
$ws.Range("P180").Value = 4.2
$ws.Range("Q180").Value = 1.727
$ws.Range("R180").Value = 0.5
$ws.Range("S180").Value = 2.05
$ws.Range("T180").Value = 1.75
$ws.Range("U180").Value = 3.5
$ws.Range("V180").Value = 1.775
$ws.Range("W180").Value = 2.025
$ws.Range("X180").Value = 2.4
$ws.Range("AA180").Value = 1.05
$ws.Range("AC180").Value = 0.7749999999999999

# Row 181
$ws.Range("B181").Value = 8271463
$ws.Range("F181").Value = "FSV GrunWeiss Ilsenburg"
$ws.Range("G181").Value = 5
$ws.Range("L181").Value = 1.166
$ws.Range("M181").Value = 7
$ws.Range("N181").Value = 9
$ws.Range("O181").Value = 1.181
$ws.Range("P181").Value = 6.5
$ws.Range("Q181").Value = 9
$ws.Range("R181").Value = -2.5
$ws.Range("S181").Value = 1.975
$ws.Range("T181").Value = 1.825
$ws.Range("U181").Value = 4.5
$ws.Range("V181").Value = 1.9
$ws.Range("W181").Value = 1.9
$ws.Range("X181").Value = 0.181
$ws.Range("AA181").Value = 0.9750000000000001
$ws.Range("AC181").Value = 0.8999999999999999

# Row 184
$ws.Range("B184").Value = 8275231
$ws.Range("F184").Value = "Sportfreunde Seligenstadt"
$ws.Range("G184").Value = 3
$ws.Range("K184").Value = "H"
$ws.Range("L184").Value = 2.4
$ws.Range("N184").Value = 2.15
$ws.Range("O184").Value = 2.875
$ws.Range("Q184").Value = 2
$ws.Range("S184").Value = 2
$ws.Range("T184").Value = 1.8
$ws.Range("U184").Value = 4
$ws.Range("V184").Value = 1.825
$ws.Range("W184").Value = 1.975
$ws.Range("X184").Value = 1.875
$ws.Range("Z184").Value = -1
$ws.Range("AA184").Value = 1
$ws.Range("AB184").Value = -1
$ws.Range("AC184").Value = 0.825
$ws.Range("AD184").Value = -1

# Row 185
$ws.Range("B185").Value = 8275229
$ws.Range("F185").Value = "SV Pars NeuIsenburg"
$ws.Range("G185").Value = 1
$ws.Range("K185").Value = "A"
$ws.Range("L185").Value = 2.5
$ws.Range("N185").Value = 2.1
$ws.Range("O185").Value = 2.5
$ws.Range("Q185").Value = 2.1
$ws.Range("S185").Value = 1.85
$ws.Range("T185").Value = 1.95
$ws.Range("U185").Value = 4.25
$ws.Range("V185").Value = 1.975
$ws.Range("W185").Value = 1.825
$ws.Range("X185").Value = -1
$ws.Range("Z185").Value = 1.1
$ws.Range("AA185").Value = -1
$ws.Range("AB185").Value = 0.95
$ws.Range("AC185").Value = -1
$ws.Range("AD185").Value = 0.825
